$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E11").Value = "La carpeta Vistas de la base de datos, buscar con el nombre V_SALDOS"
$ws.Range("E12").Value = "La carpeta Vistas de la base de datos, buscar con el nombre V_TARJETA_MENSUAL"
$ws.Range("E13").Value = "La carpeta Vistas de la base de datos, buscar con el nombre V_PAGOS_PENDIENTES"
$ws.Range("E17").Value = "La carpeta Paquetes de la base de datos, buscar con el nombre PK_GESTION_CUENTAS"
$ws.Range("E16").Value = "La carpeta Paquetes de la base de datos, buscar con el nombre PK_GESTION_CLIENTES"

$ws.Range("C60").Select()
